# Horizontal/vertical alignment fixture update:
#   - Row 4 (A4) is removed.
#   - Two new columns (D, E) are added for rows 1-3, exercising
#     vertical-only alignment styles:
#       row 1 -> vertical: top
#       row 2 -> vertical: center
#       row 3 -> vertical: general (default / no explicit alignment)
#   - Selection moves to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel XlVAlign enum values
$xlVAlignTop    = -4160
$xlVAlignCenter = -4108

# Row 1: vertical="top" on D1 (blank) and E1 ("abc")
$ws.Range("D1").VerticalAlignment = $xlVAlignTop
$ws.Range("E1").VerticalAlignment = $xlVAlignTop
$ws.Range("E1").Value = "abc"

# Row 2: vertical="center" on D2 (blank) and E2 ("abc")
$ws.Range("D2").VerticalAlignment = $xlVAlignCenter
$ws.Range("E2").VerticalAlignment = $xlVAlignCenter
$ws.Range("E2").Value = "abc"

# Row 3: D3 gets an (empty/general) alignment style applied, E3 stays
# on the default style while still holding the shared "abc" string.
$ws.Range("D3").WrapText = $false
$ws.Range("E3").Value = "abc"

# Old row 4 (A4) is removed entirely.
$ws.Range("A4").ClearContents()

# Move the active selection to D1, matching the saved selection state.
$ws.Range("D1").Select()
